$d = $word.ActiveDocument

# Locate the "Edison Achalma" byline paragraph (style "Author") that sits
# right under the "Editar: Editar" title heading.
$target = $null
foreach ($p in $d.Paragraphs) {
    $paraText = $p.Range.Text.Trim()
    if ($paraText -eq "Edison Achalma" -and $p.Style.NameLocal -eq "Author") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Create a new empty paragraph right after $target by inserting it
    # before the paragraph that currently follows $target.
    $following = $target.Next()
    $following.Range.InsertParagraphAfter()

    # The freshly created paragraph is now the one right after $target.
    $newPara = $target.Next()

    # Fill it in using raw OOXML so the result matches a normally authored
    # paragraph exactly (style + a single preserved-space run).
    $newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">Escuela Profesional de Economía, Universidad Nacional de San Cristóbal de Huamanga</w:t></w:r></w:p>'
    [void]$newPara.Range.InsertXML($newParaXml)
}
